# Update crypto price/volume table to reflect the Jan 22 2023 17:47 UTC snapshot
# (symbol list shifted by one row for B8:D17 plus refreshed price/volume figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold numeric-looking / percent-looking
# text values (e.g. "307.03", "-0.78%") that must stay literal text - format
# the whole data range as Text first so assigning them does not coerce them
# into actual numbers and lose their exact printed form.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "307.03"
$ws.Range("E2").Value = "1.05%"
$ws.Range("D3").Value = "37.05"
$ws.Range("E3").Value = "3.69%"
$ws.Range("D4").Value = "5.042"
$ws.Range("E4").Value = "-0.78%"
$ws.Range("D5").Value = "0.07902"
$ws.Range("E5").Value = "0.74%"
$ws.Range("D6").Value = "2.197"
$ws.Range("E6").Value = "-3.07%"
$ws.Range("D7").Value = "8.029"
$ws.Range("E7").Value = "-1.43%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.034"
$ws.Range("E8").Value = "0.57%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9296"
$ws.Range("E9").Value = "0.02%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.09926"
$ws.Range("E10").Value = "-0.53%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1887"
$ws.Range("E11").Value = "3.33%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08695"
$ws.Range("E12").Value = "-0.16%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03603"
$ws.Range("E13").Value = "6.40%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09958"
$ws.Range("E14").Value = "0.52%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001488"
$ws.Range("E15").Value = "0.34%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005685"
$ws.Range("E16").Value = "0.94%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.466"
$ws.Range("E17").Value = "-0.42%"
$ws.Range("E18").Value = "11.62%"
$ws.Range("D19").Value = "0.3435"
$ws.Range("E19").Value = "0.07%"
$ws.Range("D20").Value = "0.1308"
$ws.Range("E20").Value = "-0.79%"
$ws.Range("D21").Value = "4.935"
$ws.Range("E21").Value = "8.49%"
$ws.Range("E22").Value = "-1.81%"
$ws.Range("D23").Value = "0.04581"
$ws.Range("D24").Value = "0.005253"
$ws.Range("E24").Value = "16.91%"
$ws.Range("D25").Value = "0.001251"
$ws.Range("E25").Value = "1.21%"
$ws.Range("D26").Value = "0.0001399"
$ws.Range("E26").Value = "7.65%"
$ws.Range("D27").Value = "0.0002714"
$ws.Range("E27").Value = "0.73%"
$ws.Range("D39").Value = "0.01833"
$ws.Range("E39").Value = "3.65%"
$ws.Range("D40").Value = "0.04782"
$ws.Range("E40").Value = "1.70%"
$ws.Range("D41").Value = "0.007898"
$ws.Range("E41").Value = "0.15%"
$ws.Range("D42").Value = "0.1418"
$ws.Range("E42").Value = "-0.02%"
$ws.Range("D43").Value = "0.007561"
$ws.Range("E43").Value = "-10.50%"
$ws.Range("D44").Value = "0.002189"
$ws.Range("E44").Value = "-0.99%"
$ws.Range("D45").Value = "0.01056"
$ws.Range("E45").Value = "14.60%"
$ws.Range("D46").Value = "0.00006293"
$ws.Range("E46").Value = "4.04%"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").Value = "0.0005795"
$ws.Range("E48").Value = "-0.11%"
$ws.Range("D49").Value = "35.55"
$ws.Range("E49").Value = "514.84%"
$ws.Range("D50").Value = "0.002686"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").Value = "0.00002098"
$ws.Range("E51").Value = "-0.07%"
